$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 566-604: newly scraped FII/DII buy-zone data appended below the
# existing table (previously ending at row 565 / 30/05/2024).

$ws.Cells.Item(566, 1).Value = "Buying Opportunity"
$ws.Cells.Item(566, 2).Value = "support Zone"
$ws.Cells.Item(566, 3).Value = "long buildup"
$ws.Cells.Item(566, 4).Value = "Short buildup"
$ws.Cells.Item(566, 5).Value = "FII ENTERING"

$ws.Cells.Item(567, 1).Value = "ADANIPOWER"
$ws.Cells.Item(567, 2).Value = "ALPA"
$ws.Cells.Item(567, 6).Value = 755.8
$ws.Cells.Item(567, 7).Value = 81.15000000000001

$ws.Cells.Item(568, 1).Value = "ASHOKAMET"
$ws.Cells.Item(568, 2).Value = "APOLSINHOT"
$ws.Cells.Item(568, 6).Value = 24.35
$ws.Cells.Item(568, 7).Value = 1597.4

$ws.Cells.Item(569, 1).Value = "AWL"
$ws.Cells.Item(569, 2).Value = "BALKRISHNA"
$ws.Cells.Item(569, 6).Value = 355.75
$ws.Cells.Item(569, 7).Value = 25.1

$ws.Cells.Item(570, 1).Value = "DCMSHRIRAM"
$ws.Cells.Item(570, 2).Value = "BIOFILCHEM"
$ws.Cells.Item(570, 6).Value = 998.2
$ws.Cells.Item(570, 7).Value = 62.55

$ws.Cells.Item(571, 2).Value = "BSLSENETFG"
$ws.Cells.Item(571, 7).Value = 72.65000000000001

$ws.Cells.Item(572, 2).Value = "CHEMCON"
$ws.Cells.Item(572, 7).Value = 237.35

$ws.Cells.Item(573, 2).Value = "HYBRIDFIN"
$ws.Cells.Item(573, 7).Value = 12

$ws.Cells.Item(574, 2).Value = "ICICIPRULI"
$ws.Cells.Item(574, 7).Value = 545.45

$ws.Cells.Item(575, 2).Value = "INDIAGLYCO"
$ws.Cells.Item(575, 7).Value = 752.5

$ws.Cells.Item(576, 2).Value = "INSECTICID"
$ws.Cells.Item(576, 7).Value = 528.95

$ws.Cells.Item(577, 2).Value = "MTARTECH"
$ws.Cells.Item(577, 7).Value = 1800.6

$ws.Cells.Item(578, 2).Value = "NIPPOBATRY"
$ws.Cells.Item(578, 7).Value = 586.95

$ws.Cells.Item(579, 2).Value = "OSWALAGRO"
$ws.Cells.Item(579, 7).Value = 42.35

$ws.Cells.Item(580, 2).Value = "PANAMAPET"
$ws.Cells.Item(580, 7).Value = 358.9

$ws.Cells.Item(581, 2).Value = "PILITA"
$ws.Cells.Item(581, 7).Value = 12.3

$ws.Cells.Item(582, 1).NumberFormat = "@"
$ws.Cells.Item(582, 1).Value = "01/06/2024"

$ws.Cells.Item(583, 1).Value = "Buying Opportunity"
$ws.Cells.Item(583, 2).Value = "support Zone"
$ws.Cells.Item(583, 3).Value = "long buildup"
$ws.Cells.Item(583, 4).Value = "Short buildup"
$ws.Cells.Item(583, 5).Value = "FII ENTERING"

$ws.Cells.Item(584, 1).Value = "ADANIENSOL"
$ws.Cells.Item(584, 2).Value = "APLLTD"
$ws.Cells.Item(584, 6).Value = 1221.95
$ws.Cells.Item(584, 7).Value = 917.35

$ws.Cells.Item(585, 1).Value = "ADANIPOWER"
$ws.Cells.Item(585, 2).Value = "CAMPUS"
$ws.Cells.Item(585, 6).Value = 874.5
$ws.Cells.Item(585, 7).Value = 275.25

$ws.Cells.Item(586, 1).Value = "ATGL"
$ws.Cells.Item(586, 2).Value = "ENIL"
$ws.Cells.Item(586, 6).Value = 1119.4
$ws.Cells.Item(586, 7).Value = 218.15

$ws.Cells.Item(587, 1).Value = "AWL"
$ws.Cells.Item(587, 2).Value = "GOKULAGRO"
$ws.Cells.Item(587, 6).Value = 368.25
$ws.Cells.Item(587, 7).Value = 148.75

$ws.Cells.Item(588, 1).Value = "BANKBEES"
$ws.Cells.Item(588, 2).Value = "MAGNUM"
$ws.Cells.Item(588, 6).Value = 522.61
$ws.Cells.Item(588, 7).Value = 49.15

$ws.Cells.Item(589, 1).Value = "BANKETF"
$ws.Cells.Item(589, 6).Value = 512.22

$ws.Cells.Item(590, 1).Value = "EBBETF0431"
$ws.Cells.Item(590, 6).Value = 1222.46

$ws.Cells.Item(591, 1).Value = "GUJRAFFIA"
$ws.Cells.Item(591, 6).Value = 42.85

$ws.Cells.Item(592, 1).Value = "KIMS"
$ws.Cells.Item(592, 6).Value = 1832.9

$ws.Cells.Item(593, 1).Value = "LUMAXIND"
$ws.Cells.Item(593, 6).Value = 2585.55

$ws.Cells.Item(594, 1).Value = "NRL"
$ws.Cells.Item(594, 6).Value = 84.25

$ws.Cells.Item(595, 1).Value = "PSPPROJECT"
$ws.Cells.Item(595, 6).Value = 671.3

$ws.Cells.Item(596, 1).Value = "PTCIL"
$ws.Cells.Item(596, 6).Value = 9878.9

$ws.Cells.Item(597, 1).NumberFormat = "@"
$ws.Cells.Item(597, 1).Value = "03/06/2024"

$ws.Cells.Item(598, 1).Value = "Buying Opportunity"
$ws.Cells.Item(598, 2).Value = "support Zone"
$ws.Cells.Item(598, 3).Value = "long buildup"
$ws.Cells.Item(598, 4).Value = "Short buildup"
$ws.Cells.Item(598, 5).Value = "FII ENTERING"

$ws.Cells.Item(599, 1).Value = "HDFCNIF100"
$ws.Cells.Item(599, 2).Value = "BALAXI"
$ws.Cells.Item(599, 6).Value = 24.87
$ws.Cells.Item(599, 7).Value = 97.90000000000001

$ws.Cells.Item(600, 1).Value = "MID150BEES"
$ws.Cells.Item(600, 2).Value = "DCAL"
$ws.Cells.Item(600, 6).Value = 207.35
$ws.Cells.Item(600, 7).Value = 136.5

$ws.Cells.Item(601, 1).Value = "MONIFTY500"
$ws.Cells.Item(601, 2).Value = "NIRAJ"
$ws.Cells.Item(601, 6).Value = 22
$ws.Cells.Item(601, 7).Value = 42.55

$ws.Cells.Item(602, 1).Value = "NAVINIFTY"
$ws.Cells.Item(602, 6).Value = 265.03

$ws.Cells.Item(603, 1).Value = "QNIFTY"
$ws.Cells.Item(603, 6).Value = 2494.24

$ws.Cells.Item(604, 1).NumberFormat = "@"
$ws.Cells.Item(604, 1).Value = "04/06/2024"
